$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before EG (shifts EG..EH right to become EH..EI).
$ws.Columns("EG:EG").Insert()

# Row 1 header: new timestamp column.
$ws.Range("EG1").Value = "2026-02-03 06:36:33"

# Data rows: column EG mirrors the last known price in EF for rows that
# already had a numeric value there; rows with an empty EF stay empty.
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $efRef = "EF" + $r
    $efVal = $ws.Range($efRef).Value()
    if ($efVal.GetType().Name -ne "String") {
        $egRef = "EG" + $r
        $ws.Range($egRef).Value = $efVal
    }
}
